# ---------------------------------------------------------------------------
# Nexial "unitTest_numberCommand.xlsx" maintenance edit
#
# 1) [web commands]    NEW `dragTo(fromLocator,xOffset,yOffset)` - inserted
#    alphabetically into the '#system' sheet's "web" list (column V),
#    pushing every subsequent entry down by one row and extending the
#    `web` defined name from $V$2:$V$119 to $V$2:$V$120.
#
# 2) [xml commands]    NEW `beautify(xml,var)` and `minify(xml,var)` -
#    inserted alphabetically into the '#system' sheet's "xml" list
#    (column AA), pushing the trailing entries down two rows and
#    extending the `xml` defined name from $AA$2:$AA$11 to $AA$2:$AA$13.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------------
# Full, alphabetically-ordered "after" contents of column V (web commands)
# and column AA (xml commands). Re-writing the whole list (instead of doing
# a row "insert/shift") keeps every other column on the sheet untouched.
# ---------------------------------------------------------------------------
$vValues = @(
    "assertAndClick(locator,label)",
    "assertAttribute(locator,attrName,value)",
    "assertAttributeContains(locator,attrName,contains)",
    "assertAttributeNotContains(locator,attrName,contains)",
    "assertAttributeNotPresent(locator,attrName)",
    "assertAttributePresent(locator,attrName)",
    "assertChecked(locator)",
    "assertContainCount(locator,text,count)",
    "assertCssNotPresent(locator,property)",
    "assertCssPresent(locator,property,value)",
    "assertElementByAttributes(nameValues)",
    "assertElementByText(locator,text)",
    "assertElementCount(locator,count)",
    "assertElementNotPresent(locator)",
    "assertElementPresent(locator)",
    "assertFocus(locator)",
    "assertFrameCount(count)",
    "assertFramePresent(frameName)",
    "assertIECompatMode()",
    "assertIENavtiveMode()",
    "assertLinkByLabel(label)",
    "assertNotChecked(locator)",
    "assertNotFocus(locator)",
    "assertNotText(locator,text)",
    "assertNotVisible(locator)",
    "assertOneMatch(locator)",
    "assertScrollbarHNotPresent(locator)",
    "assertScrollbarHPresent(locator)",
    "assertScrollbarVNotPresent(locator)",
    "assertScrollbarVPresent(locator)",
    "assertTable(locator,row,column,text)",
    "assertText(locator,text)",
    "assertTextContains(locator,text)",
    "assertTextCount(locator,text,count)",
    "assertTextList(locator,list,ignoreOrder)",
    "assertTextMatches(text,minMatch,scrollTo)",
    "assertTextNotPresent(text)",
    "assertTextOrder(locator,descending)",
    "assertTextPresent(text)",
    "assertTitle(text)",
    "assertValue(locator,value)",
    "assertValueOrder(locator,descending)",
    "assertVisible(locator)",
    "checkAll(locator)",
    "clearLocalStorage()",
    "click(locator)",
    "clickAndWait(locator,waitMs)",
    "clickByLabel(label)",
    "clickByLabelAndWait(label,waitMs)",
    "clickWithKeys(locator,keys)",
    "close()",
    "closeAll()",
    "deselect(locator,text)",
    "deselectMulti(locator,array)",
    "dismissInvalidCert()",
    "dismissInvalidCertPopup()",
    "doubleClick(locator)",
    "doubleClickAndWait(locator,waitMs)",
    "doubleClickByLabel(label)",
    "doubleClickByLabelAndWait(label,waitMs)",
    "dragAndDrop(fromLocator,toLocator)",
    "dragTo(fromLocator,xOffset,yOffset)",
    "editLocalStorage(key,value)",
    "executeScript(var,script)",
    "focus(locator)",
    "goBack()",
    "goBackAndWait()",
    "maximizeWindow()",
    "mouseOver(locator)",
    "open(url)",
    "openAndWait(url,waitMs)",
    "openHttpBasic(url,username,password)",
    "openIgnoreTimeout(url)",
    "refresh()",
    "refreshAndWait()",
    "resizeWindow(width,height)",
    "saveAllWindowIds(var)",
    "saveAllWindowNames(var)",
    "saveAttribute(var,locator,attrName)",
    "saveCount(var,locator)",
    "saveDivsAsCsv(headers,rows,cells,nextPage,file)",
    "saveElement(var,locator)",
    "saveElements(var,locator)",
    "saveLocalStorage(var,key)",
    "saveLocation(var)",
    "savePageAs(var,sessionIdName,url)",
    "savePageAsFile(sessionIdName,url,file)",
    "saveTableAsCsv(locator,nextPageLocator,file)",
    "saveText(var,locator)",
    "saveTextArray(var,locator)",
    "saveTextSubstringAfter(var,locator,delim)",
    "saveTextSubstringBefore(var,locator,delim)",
    "saveTextSubstringBetween(var,locator,start,end)",
    "saveValue(var,locator)",
    "scrollLeft(locator,pixel)",
    "scrollRight(locator,pixel)",
    "scrollTo(locator)",
    "select(locator,text)",
    "selectFrame(locator)",
    "selectMulti(locator,array)",
    "selectMultiOptions(locator)",
    "selectText(locator)",
    "selectWindow(winId)",
    "selectWindowAndWait(winId,waitMs)",
    "selectWindowByIndex(index)",
    "selectWindowByIndexAndWait(index,waitMs)",
    "toggleSelections(locator)",
    "type(locator,value)",
    "typeKeys(locator,value)",
    "uncheckAll(locator)",
    "unselectAllText()",
    "upload(fieldLocator,file)",
    "verifyContainText(locator,text)",
    "verifyText(locator,text)",
    "wait(waitMs)",
    "waitForElementPresent(locator)",
    "waitForPopUp(winId,waitMs)",
    "waitForTextPresent(text)",
    "waitForTitle(text)"
)

$aaValues = @(
    "assertCorrectness(xml,schema)",
    "assertElementCount(xml,xpath,count)",
    "assertElementNotPresent(xml,xpath)",
    "assertElementPresent(xml,xpath)",
    "assertValue(xml,xpath,expected)",
    "assertValues(xml,xpath,array,exactOrder)",
    "assertWellformed(xml)",
    "beautify(xml,var)",
    "minify(xml,var)",
    "storeCount(xml,xpath,var)",
    "storeValue(xml,xpath,var)",
    "storeValues(xml,xpath,var)"
)

# Write column V (web), starting at row 2, 1 entry per row.
for ($i = 0; $i -lt $vValues.Length; $i++) {
    $rowNum = 2 + $i
    $addr = "V" + $rowNum
    $ws.Range($addr).Value = $vValues[$i]
}

# Write column AA (xml), starting at row 2, 1 entry per row.
for ($i = 0; $i -lt $aaValues.Length; $i++) {
    $rowNum = 2 + $i
    $addr = "AA" + $rowNum
    $ws.Range($addr).Value = $aaValues[$i]
}

# ---------------------------------------------------------------------------
# Extend the defined-name ranges to cover the newly added rows.
#   web: '#system'!$V$2:$V$119 -> '#system'!$V$2:$V$120
#   xml: '#system'!$AA$2:$AA$11 -> '#system'!$AA$2:$AA$13
# ---------------------------------------------------------------------------
$lastVRow = 1 + $vValues.Length
$lastAARow = 1 + $aaValues.Length

foreach ($n in $wb.Names) {
    if ($n.Name() -eq "web") {
        $n.RefersTo = "='#system'!`$V`$2:`$V`$" + $lastVRow
    }
    if ($n.Name() -eq "xml") {
        $n.RefersTo = "='#system'!`$AA`$2:`$AA`$" + $lastAARow
    }
}

Write-Output ("web -> V2:V" + $lastVRow)
Write-Output ("xml -> AA2:AA" + $lastAARow)
